# Apply the "env setup" update to the active workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The text values stored in B2 and B4 had their word order swapped:
#   "Test723 - LIVEHTA_723 - 1/13/2023" -> "LIVEHTA_723 - Test723 - 1/13/2023"
$ws.Range("B2").Value = "LIVEHTA_723 - Test723 - 1/13/2023"
$ws.Range("B4").Value = "LIVEHTA_723 - Test723 - 1/13/2023"

# The active selection on the sheet moved from A4 to B4.
$ws.Range("B4").Select()
